$d = $word.ActiveDocument

# "Check the validation and Test the Entire SportyShoes  Online  Application."
# becomes
# "Check the validation and Test the Entire Kitchen Story  Application."
#
# Replace "SportyShoes" (plus the trailing space that followed it) with "Kitchen",
# keeping it part of the same run as the preceding text.
$d.Content.Find.Execute(
    "Test the Entire SportyShoes ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Test the Entire Kitchen ", 2) | Out-Null

# Replace "Online" with "Story" in the following run (which keeps the double space
# before "Application").
$d.Content.Find.Execute(
    "Online  Application", $true, $false, $false, $false, $false,
    $true, 1, $false, "Story  Application", 2) | Out-Null
